$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D holds numeric-looking price text that must remain text (matches source inlineStr cells).
# Pre-format the affected numeric-looking cells as Text so assigning the new value does not
# auto-convert them into Number cells (mirrors the original non-Excel-authored inline strings).
$ws.Range("D5,D6,D7,D9,D10,D13,D15,D18,D19,D20,D21,D23,D27,D28,D29,D30,D32,D33,D34,D35,D36,D37,D40,D43,D44,D46,D47,D48").NumberFormat = "@"

$ws.Range("D2").Value = '63.815.35'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '2.632.69'
$ws.Range("E3").Value = '  +0.40%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '578.28'
$ws.Range("E5").Value = '  +0.62%  '
$ws.Range("D6").Value = '156.75'
$ws.Range("E6").Value = '  +1.00%  '
$ws.Range("D7").Value = '0.630'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = '0.118'
$ws.Range("E9").Value = '  -1.95%  '
$ws.Range("D10").Value = '5.82'
$ws.Range("E10").Value = '  +0.41%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("E12").Value = '  +0.77%  '
$ws.Range("D13").Value = '28.75'
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D14").Value = '3.109.49'
$ws.Range("E14").Value = '  +0.49%  '
$ws.Range("D15").Value = '0.0000184'
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").Value = '63.744.08'
$ws.Range("E16").Value = '  +0.17%  '
$ws.Range("D17").Value = '2.641.67'
$ws.Range("E17").Value = '  +0.95%  '
$ws.Range("D18").Value = '12.18'
$ws.Range("E18").Value = '  +0.87%  '
$ws.Range("D19").Value = '7.72'
$ws.Range("E19").Value = '  +2.76%  '
$ws.Range("D20").Value = '4.52'
$ws.Range("E20").Value = '  -2.47%  '
$ws.Range("D21").Value = '343.52'
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '68.36'
$ws.Range("E23").Value = '  +1.86%  '
$ws.Range("E24").Value = '  +9.15%  '
$ws.Range("E25").Value = '  +3.38%  '
$ws.Range("E26").Value = '  +3.62%  '
$ws.Range("D27").Value = '9.24'
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("D28").Value = '580.92'
$ws.Range("E28").Value = '  -0.51%  '
$ws.Range("D29").Value = '8.29'
$ws.Range("E29").Value = '  +4.80%  '
$ws.Range("D30").Value = '1.01'
$ws.Range("E30").Value = '  +1.46%  '
$ws.Range("E31").Value = '  -0.15%  '
$ws.Range("D32").Value = '2.05'
$ws.Range("E32").Value = '  -0.63%  '
$ws.Range("D33").Value = '1.74'
$ws.Range("E33").Value = '  +2.09%  '
$ws.Range("D34").Value = '6.67'
$ws.Range("E34").Value = '  +2.39%  '
$ws.Range("D35").Value = '5.48'
$ws.Range("E35").Value = '  +3.12%  '
$ws.Range("D36").Value = '0.404'
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("D37").Value = '19.79'
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("E39").Value = '  +2.82%  '
$ws.Range("D40").Value = '154.26'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("E41").Value = '  +8.44%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '163.35'
$ws.Range("E43").Value = '  +4.26%  '
$ws.Range("D44").Value = '24.12'
$ws.Range("E44").Value = '  +5.40%  '
$ws.Range("E45").Value = '  -0.81%  '
$ws.Range("D46").Value = '0.0588'
$ws.Range("E46").Value = '  -1.21%  '
$ws.Range("D47").Value = '0.635'
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("D48").Value = '0.100'
$ws.Range("E48").Value = '  -1.11%  '
$ws.Range("E49").Value = '  -1.15%  '
$ws.Range("D50").Value = '0.0₆0239'
$ws.Range("E50").Value = '  +2.69%  '
$ws.Range("E51").Value = '  +2.13%  '
